$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.189.02'
$ws.Range('E2').Value = '  +0.09%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.678.91'
$ws.Range('E3').Value = '  +0.54%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.41%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.67'
$ws.Range('E5').Value = '  -0.37%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5259'
$ws.Range('E6').Value = '  +3.28%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.003'
$ws.Range('E7').Value = '  -0.41%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2700'
$ws.Range('E8').Value = '  +1.88%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06416'
$ws.Range('E9').Value = '  +1.71%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.01'
$ws.Range('E10').Value = '  +0.98%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07442'
$ws.Range('E11').Value = '  +0.87%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.696.97'
$ws.Range('E12').Value = '  +1.44%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.532'
$ws.Range('E13').Value = '  -0.17%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5856'
$ws.Range('E14').Value = '  +1.79%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008550'
$ws.Range('E15').Value = '  +0.61%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.51'
$ws.Range('E16').Value = '  -0.55%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.217.89'
$ws.Range('E17').Value = '  -0.23%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.954'
$ws.Range('E18').Value = '  -0.73%  '

$ws.Range('E19').Value = '  -0.22%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.82'
$ws.Range('E20').Value = '  -0.21%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '190.87'
$ws.Range('E21').Value = '  +2.28%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.224'
$ws.Range('E22').Value = '  +0.42%  '

$ws.Range('E23').Value = '  -0.35%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '145.26'
$ws.Range('E24').Value = '  +1.28%  '

$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1253'
$ws.Range('E25').Value = '  +7.19%  '

$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.645'
$ws.Range('E26').Value = '  +1.06%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.81'
$ws.Range('E27').Value = '  +0.84%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06592'
$ws.Range('E28').Value = '  +14.46%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.334'
$ws.Range('E29').Value = '  +1.37%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.319'
$ws.Range('E30').Value = '  -0.87%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.598'
$ws.Range('E31').Value = '  +2.67%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.543'
$ws.Range('E32').Value = '  +1.40%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.673'

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.021'
$ws.Range('E34').Value = '  +1.99%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6183'
$ws.Range('E35').Value = '  +3.44%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.367'
$ws.Range('E36').Value = '  -0.13%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.693'
$ws.Range('E37').Value = '  +2.05%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.280'
$ws.Range('E38').Value = '  +6.19%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.099.01'
$ws.Range('E39').Value = '  +0.03%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01601'
$ws.Range('E40').Value = '  +0.06%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8725'
$ws.Range('E41').Value = '  +1.02%  '

$ws.Range('E42').Value = '  +0.59%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.00'
$ws.Range('E43').Value = '  +1.62%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.820.88'
$ws.Range('E44').Value = '  -0.06%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000112'
$ws.Range('E45').Value = '  -1.51%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.69'
$ws.Range('E46').Value = '  +0.93%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.144'
$ws.Range('E47').Value = '  +1.23%  '

$ws.Range('E48').Value = '  -0.26%  '

$ws.Range('E49').Value = '  +0.71%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4284'
$ws.Range('E50').Value = '  -0.76%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.013'
$ws.Range('E51').Value = '  +3.56%  '

